$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("E2").Value = '  -0.20%  '
Set-TextValue $ws.Range("D2") '43.684.78'

$ws.Range("E3").Value = '  -0.45%  '
Set-TextValue $ws.Range("D3") '2.281.24'

$ws.Range("E5").Value = '  +10.45%  '
Set-TextValue $ws.Range("D5") '113.87'

$ws.Range("E6").Value = '  -1.26%  '
Set-TextValue $ws.Range("D6") '267.10'

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  +0.27%  '

$ws.Range("E9").Value = '  +0.56%  '
Set-TextValue $ws.Range("D9") '0.611'

$ws.Range("E10").Value = '  +5.68%  '

$ws.Range("E11").Value = '  +0.25%  '

$ws.Range("E12").Value = '  +8.69%  '
Set-TextValue $ws.Range("D12") '8.83'

$ws.Range("E13").Value = '  +0.55%  '
Set-TextValue $ws.Range("D13") '0.108'

$ws.Range("E14").Value = '  +0.88%  '
Set-TextValue $ws.Range("D14") '15.72'

$ws.Range("E15").Value = '  -0.39%  '
Set-TextValue $ws.Range("D15") '2.623.88'

$ws.Range("E16").Value = '  +2.01%  '
Set-TextValue $ws.Range("D16") '0.873'

$ws.Range("E17").Value = '  -0.38%  '
Set-TextValue $ws.Range("D17") '2.278.91'

$ws.Range("E18").Value = '  -0.66%  '
Set-TextValue $ws.Range("D18") '43.467.63'

$ws.Range("E19").Value = '  -1.36%  '

$ws.Range("E20").Value = '  +12.70%  '
Set-TextValue $ws.Range("D20") '7.07'

$ws.Range("E21").Value = '  -0.46%  '
Set-TextValue $ws.Range("D21") '71.97'

$ws.Range("E22").Value = '  -3.75%  '
Set-TextValue $ws.Range("D22") '2.41'

$ws.Range("E23").Value = '  +7.09%  '
Set-TextValue $ws.Range("D23") '9.93'

$ws.Range("E24").Value = '  -0.51%  '
Set-TextValue $ws.Range("D24") '232.37'

$ws.Range("E25").Value = '  -0.39%  '
Set-TextValue $ws.Range("D25") '2.87'

$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("E27").Value = '  +2.08%  '
Set-TextValue $ws.Range("D27") '11.53'

$ws.Range("E28").Value = '  +0.49%  '
Set-TextValue $ws.Range("D28") '40.83'

$ws.Range("E29").Value = '  -1.52%  '

$ws.Range("E30").Value = '  -1.22%  '

Set-TextValue $ws.Range("D31") '173.37'

$ws.Range("E32").Value = '  -1.58%  '
Set-TextValue $ws.Range("D32") '21.49'

$ws.Range("E33").Value = '  +0.69%  '
Set-TextValue $ws.Range("D33") '0.0909'

$ws.Range("E34").Value = '  +2.38%  '
Set-TextValue $ws.Range("D34") '5.65'

$ws.Range("E35").Value = '  +0.24%  '

$ws.Range("E36").Value = '  -5.79%  '
Set-TextValue $ws.Range("D36") '4.61'

$ws.Range("E37").Value = '  -2.33%  '
Set-TextValue $ws.Range("D37") '0.0351'

$ws.Range("E38").Value = '  -4.48%  '

$ws.Range("E39").Value = '  +5.80%  '
Set-TextValue $ws.Range("D39") '3.74'

$ws.Range("E40").Value = '  +18.61%  '
Set-TextValue $ws.Range("D40") '14.47'

$ws.Range("E41").Value = '  +14.42%  '
Set-TextValue $ws.Range("D41") '75.09'

$ws.Range("E42").Value = '  +3.69%  '
Set-TextValue $ws.Range("D42") '2.42'

$ws.Range("E43").Value = '  +0.35%  '

$ws.Range("E44").Value = '  +17.23%  '
Set-TextValue $ws.Range("D44") '6.20'

$ws.Range("E45").Value = '  +0.13%  '

$ws.Range("E46").Value = '  +0.37%  '

$ws.Range("E47").Value = '  -1.71%  '
Set-TextValue $ws.Range("D47") '8.66'

$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E48").Value = '  +2.87%  '
Set-TextValue $ws.Range("D48") '102.06'

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("E49").Value = '  -2.23%  '
Set-TextValue $ws.Range("D49") '0.0995'

$ws.Range("B50").Value = 'TrustWalletToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("E50").Value = '  +1.76%  '
Set-TextValue $ws.Range("D50") '1.25'

$ws.Range("E51").Value = '  +2.81%  '
Set-TextValue $ws.Range("D51") '0.453'

